$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 32 for day 31 of the month
$ws.Range("A32").Value = 31

$ws.Range("B2").Value = "火"
$ws.Range("D2").Value = "清水在宅作業"
$ws.Range("E2").Value = "09:02"
$ws.Range("F2").Value = "18:04"
$ws.Range("G2").Value = "1:00"

$ws.Range("B3").Value = "水"
$ws.Range("E3").Value = "09:37"
$ws.Range("F3").Value = "18:37"

$ws.Range("B4").Value = "木"
$ws.Range("E4").Value = "09:01"
$ws.Range("F4").Value = "18:26"

$ws.Range("B5").Value = "金"
$ws.Range("E5").Value = "09:20"
$ws.Range("F5").Value = "18:29"

$ws.Range("B6").Value = "土"
$ws.Range("D6").Value = ""
$ws.Range("E6").Value = ""
$ws.Range("F6").Value = ""
$ws.Range("G6").Value = ""

$ws.Range("B7").Value = "日"
$ws.Range("D7").Value = ""
$ws.Range("E7").Value = ""
$ws.Range("F7").Value = ""
$ws.Range("G7").Value = ""

$ws.Range("B8").Value = "月"
$ws.Range("D8").Value = "清水在宅作業"
$ws.Range("E8").Value = "09:01"
$ws.Range("F8").Value = "18:40"
$ws.Range("G8").Value = "1:00"

$ws.Range("B9").Value = "火"
$ws.Range("D9").Value = "清水在宅作業"
$ws.Range("E9").Value = "09:00"
$ws.Range("F9").Value = "20:01"
$ws.Range("G9").Value = "1:00"

$ws.Range("B10").Value = "水"
$ws.Range("E10").Value = "09:06"
$ws.Range("F10").Value = "19:03"

$ws.Range("B11").Value = "木"
$ws.Range("E11").Value = "09:18"
$ws.Range("F11").Value = "18:21"

$ws.Range("B12").Value = "金"
$ws.Range("E12").Value = "08:54"
$ws.Range("F12").Value = "18:17"

$ws.Range("B13").Value = "土"
$ws.Range("D13").Value = ""
$ws.Range("E13").Value = ""
$ws.Range("F13").Value = ""
$ws.Range("G13").Value = ""

$ws.Range("B14").Value = "日"
$ws.Range("D14").Value = ""
$ws.Range("E14").Value = ""
$ws.Range("G14").Value = ""

$ws.Range("B15").Value = "月"
$ws.Range("D15").Value = "清水在宅作業"
$ws.Range("E15").Value = "09:29"
$ws.Range("F15").Value = "18:38"
$ws.Range("G15").Value = "1:00"

$ws.Range("B16").Value = "火"
$ws.Range("D16").Value = "清水在宅作業"
$ws.Range("E16").Value = "09:12"
$ws.Range("F16").Value = "18:36"
$ws.Range("G16").Value = "1:00"

$ws.Range("B17").Value = "水"
$ws.Range("E17").Value = "09:30"
$ws.Range("F17").Value = "18:30"

$ws.Range("B18").Value = "木"
$ws.Range("E18").Value = "09:41"
$ws.Range("F18").Value = "19:01"

$ws.Range("B19").Value = "金"
$ws.Range("E19").Value = "09:05"
$ws.Range("F19").Value = "18:21"

$ws.Range("B20").Value = "土"
$ws.Range("D20").Value = ""
$ws.Range("E20").Value = ""
$ws.Range("F20").Value = ""
$ws.Range("G20").Value = ""

$ws.Range("B21").Value = "日"
$ws.Range("D21").Value = ""
$ws.Range("E21").Value = ""
$ws.Range("F21").Value = ""
$ws.Range("G21").Value = ""

$ws.Range("B22").Value = "月"

$ws.Range("B23").Value = "火"
$ws.Range("D23").Value = "清水在宅作業"
$ws.Range("E23").Value = "10:31"
$ws.Range("F23").Value = "20:04"
$ws.Range("G23").Value = "1:00"

$ws.Range("B24").Value = "水"
$ws.Range("E24").Value = "10:50"
$ws.Range("F24").Value = "20:11"

$ws.Range("B25").Value = "木"
$ws.Range("E25").Value = "09:31"
$ws.Range("F25").Value = "19:10"

$ws.Range("B26").Value = "金"
$ws.Range("E26").Value = "11:01"
$ws.Range("F26").Value = "20:04"

$ws.Range("B27").Value = "土"

$ws.Range("B28").Value = "日"
$ws.Range("D28").Value = ""
$ws.Range("E28").Value = ""
$ws.Range("F28").Value = ""
$ws.Range("G28").Value = ""

$ws.Range("B29").Value = "月"
$ws.Range("D29").Value = "清水在宅作業"
$ws.Range("E29").Value = "09:24"
$ws.Range("F29").Value = "18:29"
$ws.Range("G29").Value = "1:00"

$ws.Range("B30").Value = "火"
$ws.Range("D30").Value = "清水在宅作業"
$ws.Range("E30").Value = "09:43"
$ws.Range("F30").Value = "21:00"
$ws.Range("G30").Value = "1:00"

$ws.Range("B31").Value = "水"
$ws.Range("E31").Value = "09:06"
$ws.Range("F31").Value = "18:25"

$ws.Range("B32").Value = "木"
$ws.Range("D32").Value = "清水在宅作業"
$ws.Range("E32").Value = "09:06"
$ws.Range("F32").Value = "18:18"
$ws.Range("G32").Value = "1:00"